# "error solve ifrs list"
# Replaces the stale/incorrect financial figures (rows 2-6, columns D:AJ) with
# the corrected values, and removes the erroneous extra forecast rows 7-9
# (their D:AI figures were bogus placeholders) leaving only the row
# number / "연간" / period-label columns (A/B/C) intact for those rows.
# Row 6 also drops the AG6/AH6 (현금DPS/현금배당수익률) cells entirely, since
# they no longer have reliable source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 6601
$ws.Range("E2").Value = 77
$ws.Range("F2").Value = 77
$ws.Range("G2").Value = 61
$ws.Range("H2").Value = 52
$ws.Range("I2").Value = 56
$ws.Range("J2").Value = -4
$ws.Range("K2").Value = 6728
$ws.Range("L2").Value = 1687
$ws.Range("M2").Value = 5041
$ws.Range("N2").Value = 4819
$ws.Range("O2").Value = 221
$ws.Range("P2").Value = 250
$ws.Range("Q2").Value = 250
$ws.Range("R2").Value = -28
$ws.Range("S2").Value = -250
$ws.Range("T2").Value = 66
$ws.Range("U2").Value = 184
$ws.Range("V2").Value = 731
$ws.Range("W2").Value = 1.16
$ws.Range("X2").Value = 0.78
$ws.Range("Y2").Value = 1.16
$ws.Range("Z2").Value = 0.75
$ws.Range("AA2").Value = 33.46
$ws.Range("AB2").Value = 1836.27
$ws.Range("AC2").Value = 1114
$ws.Range("AD2").Value = 25.89
$ws.Range("AE2").Value = 99799
$ws.Range("AF2").Value = 0.29
$ws.Range("AG2").Value = 350
$ws.Range("AH2").Value = 1.21
$ws.Range("AI2").Value = 30.31
$ws.Range("AJ2").Value = 5004949

# Row 3
$ws.Range("D3").Value = 6552
$ws.Range("E3").Value = 169
$ws.Range("F3").Value = 169
$ws.Range("G3").Value = 135
$ws.Range("H3").Value = 97
$ws.Range("I3").Value = 93
$ws.Range("J3").Value = 5
$ws.Range("K3").Value = 7199
$ws.Range("L3").Value = 2076
$ws.Range("M3").Value = 5124
$ws.Range("N3").Value = 4882
$ws.Range("O3").Value = 242
$ws.Range("P3").Value = 250
$ws.Range("Q3").Value = -401
$ws.Range("R3").Value = 89
$ws.Range("S3").Value = 345
$ws.Range("T3").Value = 96
$ws.Range("U3").Value = -498
$ws.Range("V3").Value = 1138
$ws.Range("W3").Value = 2.59
$ws.Range("X3").Value = 1.49
$ws.Range("Y3").Value = 1.91
$ws.Range("Z3").Value = 1.4
$ws.Range("AA3").Value = 40.51
$ws.Range("AB3").Value = 1858.61
$ws.Range("AC3").Value = 1853
$ws.Range("AD3").Value = 18.89
$ws.Range("AE3").Value = 101096
$ws.Range("AF3").Value = 0.35
$ws.Range("AG3").Value = 450
$ws.Range("AH3").Value = 1.29
$ws.Range("AI3").Value = 23.44
$ws.Range("AJ3").Value = 5004949

# Row 4
$ws.Range("D4").Value = 6573
$ws.Range("E4").Value = 244
$ws.Range("F4").Value = 249
$ws.Range("G4").Value = 215
$ws.Range("H4").Value = 164
$ws.Range("I4").Value = 173
$ws.Range("J4").Value = -9
$ws.Range("K4").Value = 7301
$ws.Range("L4").Value = 2112
$ws.Range("M4").Value = 5189
$ws.Range("N4").Value = 4959
$ws.Range("O4").Value = 230
$ws.Range("P4").Value = 250
$ws.Range("Q4").Value = 269
$ws.Range("R4").Value = -144
$ws.Range("S4").Value = -142
$ws.Range("T4").Value = 94
$ws.Range("U4").Value = 175
$ws.Range("V4").Value = 1062
$ws.Range("W4").Value = 3.71
$ws.Range("X4").Value = 2.5
$ws.Range("Y4").Value = 3.51
$ws.Range("Z4").Value = 2.27
$ws.Range("AA4").Value = 40.71
$ws.Range("AB4").Value = 1890.88
$ws.Range("AC4").Value = 3453
$ws.Range("AD4").Value = 8.92
$ws.Range("AE4").Value = 102681
$ws.Range("AF4").Value = 0.3
$ws.Range("AG4").Value = 600
$ws.Range("AH4").Value = 1.95
$ws.Range("AI4").Value = 16.77
$ws.Range("AJ4").Value = 5004949

# Row 5
$ws.Range("D5").Value = 6779
$ws.Range("E5").Value = 97
$ws.Range("F5").Value = 97
$ws.Range("G5").Value = 139
$ws.Range("H5").Value = 80
$ws.Range("I5").Value = 98
$ws.Range("J5").Value = -19
$ws.Range("K5").Value = 7714
$ws.Range("L5").Value = 2508
$ws.Range("M5").Value = 5206
$ws.Range("N5").Value = 5000
$ws.Range("O5").Value = 206
$ws.Range("P5").Value = 250
$ws.Range("Q5").Value = 15
$ws.Range("R5").Value = -265
$ws.Range("S5").Value = 332
$ws.Range("T5").Value = 335
$ws.Range("U5").Value = -321
$ws.Range("V5").Value = 1400
$ws.Range("W5").Value = 1.43
$ws.Range("X5").Value = 1.18
$ws.Range("Y5").Value = 1.98
$ws.Range("Z5").Value = 1.06
$ws.Range("AA5").Value = 48.17
$ws.Range("AB5").Value = 1917.51
$ws.Range("AC5").Value = 1968
$ws.Range("AD5").Value = 11.38
$ws.Range("AE5").Value = 103543
$ws.Range("AF5").Value = 0.22
$ws.Range("AG5").Value = 300
$ws.Range("AH5").Value = 1.34
$ws.Range("AI5").Value = 14.71
$ws.Range("AJ5").Value = 5004949

# Row 6
$ws.Range("D6").Value = 7234
$ws.Range("E6").Value = -145
$ws.Range("F6").Value = -145
$ws.Range("G6").Value = -446
$ws.Range("H6").Value = -330
$ws.Range("I6").Value = -344
$ws.Range("K6").Value = 7582
$ws.Range("L6").Value = 2717
$ws.Range("M6").Value = 4865
$ws.Range("N6").Value = 4650
$ws.Range("P6").Value = 250
$ws.Range("Q6").Value = -331
$ws.Range("R6").Value = -55
$ws.Range("S6").Value = 259
$ws.Range("T6").Value = 114
$ws.Range("U6").Value = -445
$ws.Range("V6").Value = 1724
$ws.Range("W6").Value = -2
$ws.Range("X6").Value = -4.56
$ws.Range("Y6").Value = -7.12
$ws.Range("Z6").Value = -4.31
$ws.Range("AA6").Value = 55.85
$ws.Range("AB6").Value = 1775.78
$ws.Range("AC6").Value = -6866
$ws.Range("AD6").Value = -2.37
$ws.Range("AE6").Value = 96282
$ws.Range("AF6").Value = 0.17
$ws.Range("AI6").Value = -1.41
$ws.Range("AJ6").Value = 5004949
$ws.Range("AG6:AH6").ClearContents()

# Rows 7-9: clear all data columns except A/B/C (row number, period label)
$ws.Range("D7:AI7").ClearContents()
$ws.Range("D8:AI8").ClearContents()
$ws.Range("D9:AI9").ClearContents()
